# Adicionando comentarios no codigo
# Atualiza a tabela de previsao do tempo (Manaus): os dados avancam um dia
# e uma nova linha (dia 11) e adicionada ao final da tabela.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Garante que as colunas de umidade (percentuais) sejam tratadas como texto,
# preservando valores como "61%" em vez de serem convertidos em numeros/porcentagens.
$ws.Range("D2:E11").NumberFormat = "@"

# Dados atualizados para cada dia (Dia, Temp. maxima, Temp. minima,
# Umidade do dia, Umidade da noite, Indice UV do dia, Indice UV da noite)
$data = @{
    2  = @("qua. 25", "35°", "26°", "61%", "83%", "Extremo", "0 de 11")
    3  = @("qui. 26", "35°", "26°", "61%", "82%", "Extremo", "0 de 11")
    4  = @("sex. 27", "34°", "26°", "64%", "82%", "Extremo", "0 de 11")
    5  = @("sáb. 28", "34°", "26°", "66%", "86%", "Extremo", "0 de 11")
    6  = @("dom. 29", "32°", "25°", "73%", "87%", "Extremo", "0 de 11")
    7  = @("seg. 30", "32°", "25°", "72%", "87%", "Extremo", "0 de 11")
    8  = @("ter. 01", "31°", "25°", "75%", "91%", "Extremo", "0 de 11")
    9  = @("qua. 02", "34°", "25°", "68%", "89%", "Extremo", "0 de 11")
    10 = @("qui. 03", "33°", "25°", "71%", "91%", "Extremo", "0 de 11")
    11 = @("sex. 04", "33°", "25°", "69%", "89%", "Extremo", "0 de 11")
}

foreach ($rowNum in $data.Keys) {
    $rowValues = $data[$rowNum]
    for ($col = 1; $col -le 7; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $rowValues[$col - 1]
    }
}

# A tabela agora tem uma linha extra (dia 11), atualiza a dimensao da planilha.
$ws.Range("A1:G11").Select() | Out-Null
